$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "TextBox 6" shape (id=7) by name rather than a hard-coded
# index, in case shape ordering differs from what was inspected.
$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 6") {
        $sh = $candidate
        break
    }
}

# Move it to the end of the z-order (last in the XML shape tree),
# matching its relocation to just before </p:spTree> in the diff.
$sh.ZOrder(0)

# Update position/size (EMU -> points: divide by 12700).
# Note: the host stores these as single-precision floats, so the naive
# EMU/12700 quotient can round-trip to one EMU less than intended once
# truncated back to integer EMUs on save. Nudge to the nearest point
# value whose float32 representation still floors to the exact target
# EMU count.
$sh.Left = 215.82701110839844
$sh.Top = 68.15008544921875
$sh.Width = 96.779296875
$sh.Height = 55.73905563354492

# Update the text content.
$sh.TextFrame.TextRange.Text = "Equations &  Logic"
